# Applies the "after mids syllabus added for algo" update to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C12: TPL Week 9 now has 1 uploaded lecture (was "Not uploaded").
$ws.Range("C12").Value = "Topic: name- Week 9, lectures- 1 Lecture, duration- 01:18;`nVideo: link- https://drive.google.com/file/d/13wRNzU6PSClUQ2hdxH6xlae5zPEjHR2d/preview, name- TPL Week # 9.mp4, duration- 01:18;"

# D14: SRE Week 11 now has 2 uploaded lectures (was "Not uploaded").
$ws.Range("D14").Value = "Topic: name- Week 11, lectures- 2 Lectures, duration- 00:56;`nVideo: link- https://drive.google.com/file/d/1iGjE5XenE92z_rJ_1-lJD1LYZIY4clXC/preview, name- SRE Week # 11 Part 1.mp4, duration- 00:42;`nVideo: link- https://drive.google.com/file/d/1lM0So4usTGuTnFKxlFqhn04Cv3P0bVqP/preview, name- SRE Week # 11 Part 2.mp4, duration- 00:14;"

# Row 15: new Week 12 entries for TPL / SRE (not uploaded yet) and ALGO (2 lectures uploaded).
$ws.Range("C15").Value = "Topic: name- Week 12, lectures- Not uploaded, duration- -;`nNote: heading- Lectures not uploaded yet;"
$ws.Range("D15").Value = "Topic: name- Week 12, lectures- Not uploaded, duration- -;`nNote: heading- Lectures not uploaded yet;"
$ws.Range("E15").Value = "Topic: name- Week 12, lectures- 2 Lectures, duration- 01:36;`nVideo: link- https://drive.google.com/file/d/1kOuhYo1poWvPz8T0eoBRzsXcoKoY7U12/preview, name- AD&AA Week # 12 Part 1.mp4, duration- 00:58;`nVideo: link- https://drive.google.com/file/d/1t2dXItzIgFA2pX8ZnqEudZNZBvRmACIF/preview, name- AD&AA Week # 12 Part 2.mp4, duration- 00:38;"

# Row 16: new Final Term Exam syllabus info for ALGO (after mids).
$ws.Range("E16").Value = "Topic: name- Final Term Exam, lectures- , duration- 27 Dec | 16:00-17:30;`nAssignment: name- Presentations List, link- presentation list.xlsx;`nSlides: slide- lec5.ppt;`nSlides: slide- Lec6.pptx;`nSlides: slide- lec7.pptx;`nSlides: slide- lec8.pdf;`nSlides: slide- lec9.pptx;`nSlides: slide- lec10.pdf;`nNote: heading- Syllabus not final yet. But these lectures are taught after mids;"

# Match formatting of the other populated "Topic" cells (bold, wrapped, left/middle
# aligned) by copying the format from an existing cell that already uses that style.
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15:E15").PasteSpecial(-4122)
$ws.Range("E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(15).RowHeight = 170
$ws.Rows.Item(16).RowHeight = 204

# Keep the view roughly where the author left it.
$ws.Range("E16").Select()
